# include Liyunet comments on Oromiya
# Rewrite the weredas lookup table for the Oromiya region sheet so that it
# reflects the corrected / updated wereda names and zone groupings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previously used range (header + up to 60 data rows, 3 columns)
$ws.Range("A1:C100").ClearContents()

# Header row
$ws.Range("A1").Value = "Region"
$ws.Range("B1").Value = "Zone"
$ws.Range("C1").Value = "Wereda"

$data = @(
    @("oromiya","arsi","arsi"),
    @("oromiya","arsi","deksis"),
    @("oromiya","arsi","enkelo wabe"),
    @("oromiya","arsi","hitosa"),
    @("oromiya","arsi","robe"),
    @("oromiya","arsi","sire"),
    @("oromiya","arsi","sude"),
    @("oromiya","arsi","zeway dugda"),
    @("oromiya","bale","agarfa"),
    @("oromiya","bale","berbere"),
    @("oromiya","bale","delo mena"),
    @("oromiya","bale","gasera"),
    @("oromiya","bale","sinana"),
    @("oromiya","east hararge","bedeno"),
    @("oromiya","east hararge","burqaa"),
    @("oromiya","east hararge","deder"),
    @("oromiya","east hararge","fedis"),
    @("oromiya","east hararge","g /muxii"),
    @("oromiya","east hararge","gole oda"),
    @("oromiya","east hararge","goroo  muxii"),
    @("oromiya","east hararge","gurawa"),
    @("oromiya","east hararge","gurawaa"),
    @("oromiya","east hararge","gursum"),
    @("oromiya","east hararge","haromaya"),
    @("oromiya","east hararge","jaarsoo"),
    @("oromiya","east hararge","meta"),
    @("oromiya","east hararge","midega tola"),
    @("oromiya","east hararge","oda muda"),
    @("oromiya","jima","gera"),
    @("oromiya","jima","j/z/ social affairs"),
    @("oromiya","jima","kersa"),
    @("oromiya","jima","manna"),
    @("oromiya","jima","seka chekorsa"),
    @("oromiya","west arsi","adaba"),
    @("oromiya","west arsi","aminya"),
    @("oromiya","west arsi","ar/lixa"),
    @("oromiya","west arsi","arsi negele"),
    @("oromiya","west arsi","arsii"),
    @("oromiya","west arsi","dodola"),
    @("oromiya","west arsi","hitosa"),
    @("oromiya","west arsi","kofele"),
    @("oromiya","west arsi","kokossa"),
    @("oromiya","west arsi","kore"),
    @("oromiya","west arsi","robe"),
    @("oromiya","west arsi","sire"),
    @("oromiya","west arsi","sude"),
    @("oromiya","west arsi","xana"),
    @("oromiya","west hararge","bedessa town"),
    @("oromiya","west hararge","chiro zuria"),
    @("oromiya","west hararge","daro lebu"),
    @("oromiya","west hararge","doba"),
    @("oromiya","west hararge","guba qoricha"),
    @("oromiya","west hararge","habro"),
    @("oromiya","west hararge","tulo")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $row = $row + 1
}
